# Introduction to Quantum Cryptography - small title + bibliography cleanup
$p = $ppt.ActivePresentation

# -----------------------------------------------------------------
# Slide 1: title change
#   "Quantum and Post Quantum Cryptography"
#     -> "Introduction to Quantum Cryptography"
#   split into 3 runs: "Introduction to " (en-US) + "Quantum " (el) + "Cryptography" (el)
# -----------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$titleShape = $s1.Shapes.Item(1)
$titleRange = $titleShape.TextFrame.TextRange
$titleRange.Text = "Introduction to Quantum Cryptography"

# Force the paragraph to split into the "Introduction to " / "Quantum " / "Cryptography" runs
$part1 = $titleShape.TextFrame.TextRange.Characters(1, 16)
$part1.Text = $part1.Text
$part2 = $titleShape.TextFrame.TextRange.Characters(17, 8)
$part2.Text = $part2.Text

# Mark the newly-added lead-in text as English; the rest keeps its original Greek language tag
$leadIn = $titleShape.TextFrame.TextRange.Characters(1, 16)
$leadIn.LanguageID = "en-US"

# -----------------------------------------------------------------
# Slide 23: bibliography - collapse runs that were split mid-word/url
#   back into single runs (no visible text change, only run structure)
# -----------------------------------------------------------------
$s23 = $p.Slides.Item(23)
$bibShape = $s23.Shapes.Item(2)
$bibRange = $bibShape.TextFrame.TextRange

# ", volume 175, page 8. New York, " + "1984. " -> one run
$merge1 = $bibRange.Characters(189, 38)
$merge1.Text = ", volume 175, page 8. New York, 1984. "

# "http" + "://" + "researcher.watson.ibm.com/researcher/files/us-bennetc/BB84highest.pdf" -> one run
$merge2 = $bibShape.TextFrame.TextRange.Characters(227, 77)
$merge2.Text = "http://researcher.watson.ibm.com/researcher/files/us-bennetc/BB84highest.pdf"

# "https://" + "en.wikipedia.org/wiki/Quantum_key_distribution" -> one run
$merge3 = $bibShape.TextFrame.TextRange.Characters(370, 55)
$merge3.Text = "https://en.wikipedia.org/wiki/Quantum_key_distribution"

# "https://" + "en.wikipedia.org/wiki/Schr%C3%B6dinger%27s_cat" -> one run
$merge4 = $bibShape.TextFrame.TextRange.Characters(425, 55)
$merge4.Text = "https://en.wikipedia.org/wiki/Schr%C3%B6dinger%27s_cat"
